$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the tiny floating-point correction on A4's timestamp
$ws.Range("A4").Value = 45875.12516701389

# Append the new row of data (row 5)
$ws.Range("A5").Value = 45875.20850464647
$ws.Range("A5").NumberFormat = $ws.Range("A4").NumberFormat

$ws.Range("B5").Value = 2025
$ws.Range("C5").Value = 23
$ws.Range("D5").Value = 13.67
$ws.Range("E5").Value = 92.84
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0.01
$ws.Range("H5").Value = "NW"
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = "05:00:14"
